$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.460.16'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '3.689.14'
$ws.Range('E3').Value = '  -3.06%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '682.06'
$ws.Range('E5').Value = '  -3.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.67'
$ws.Range('E6').Value = '  -4.54%  '
$ws.Range('D7').Value = '3.687.90'
$ws.Range('E7').Value = '  -3.10%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -4.30%  '
$ws.Range('E10').Value = '  -7.71%  '
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('E12').Value = '  -2.55%  '
$ws.Range('E13').Value = '  -4.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.67'
$ws.Range('E14').Value = '  -6.06%  '
$ws.Range('D15').Value = '4.312.93'
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').Value = '3.683.17'
$ws.Range('E16').Value = '  -3.22%  '
$ws.Range('D17').Value = '69.517.19'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.30'
$ws.Range('E19').Value = '  -6.19%  '
$ws.Range('E20').Value = '  -6.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '482.96'
$ws.Range('E21').Value = '  -3.28%  '
$ws.Range('E22').Value = '  -6.92%  '
$ws.Range('E23').Value = '  -7.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.32'
$ws.Range('E24').Value = '  -4.47%  '
$ws.Range('D25').Value = '3.834.89'
$ws.Range('E25').Value = '  -3.04%  '
$ws.Range('E26').Value = '  -8.51%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.48'
$ws.Range('E28').Value = '  -4.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.55'
$ws.Range('E29').Value = '  -7.15%  '
$ws.Range('E30').Value = '  -8.46%  '
$ws.Range('E31').Value = '  -10.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.89'
$ws.Range('E32').Value = '  -5.75%  '
$ws.Range('E33').Value = '  -7.39%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.16'
$ws.Range('E34').Value = '  -6.32%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.168'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').Value = '3.659.72'
$ws.Range('E37').Value = '  -2.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.51'
$ws.Range('E38').Value = '  -5.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.34'
$ws.Range('E39').Value = '  +6.71%  '
$ws.Range('E40').Value = '  -6.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.25'
$ws.Range('E41').Value = '  -4.17%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  -6.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '161.54'
$ws.Range('E45').Value = '  -3.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.36'
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('E47').Value = '  -11.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.26'
$ws.Range('E48').Value = '  +7.99%  '
$ws.Range('E49').Value = '  -8.11%  '
$ws.Range('E50').Value = '  +0.33%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.12'
$ws.Range('E51').Value = '  -1.47%  '
